$wb = $excel.ActiveWorkbook

# --- Worksheets involved ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Datetime stamps refreshed for the new handoff generation ---
$wsOverview.Range("G2").Value = "2016-08-12 06:47:49"
$wsDeDe.Range("H2").Value     = "2016-08-12 06:47:49"
$wsZhCn.Range("H2").Value     = "2016-08-12 06:47:42"

# --- Column widths widen to fit the new, longer "Ready for handoff" text ---
# (The ColumnWidth COM property here is quantized in ~1/6 character-width
# steps, so we pick the input that lands on the closest attainable width.)
$wsOverview.Range("E1").ColumnWidth = 16.333333333333332
$wsOverview.Range("F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth     = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth     = 16.333333333333332
